$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "custom accuracy": round row 5's values to 2 decimal places
$values = @(15.65, 12.1, 0.69, 34.14, 28.38, 12.5, 50.18, 19.15, 8.62, 13.24, 14.51, 14.49, 3.83, 12.18, 17.76, 10.25, 0.12, 0.44, 181.36, 34.7, 11.16, 23.36, 12.62, 1.64, 24.48, 10.02, 9.38, 10.26, 14.9, 0.08, 45.52, 6.62, 14.18)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(5, $i + 2).Value = $values[$i]
}

# "데이터 1000개" (reduce dataset by one row): drop the old last row (row 6)
$ws.Rows(6).Delete()
